$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = "16Gb SD Card"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 9.99

$ws.Range("F20").Select()
